$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '29.342.58'
$ws.Range("E2").Value = '  -0.25%  '

$ws.Range("D3").Value = '1.846.62'
$ws.Range("E3").Value = '  -0.17%  '

$ws.Range("D4").Value = '''0.9981'
$ws.Range("E4").Value = '  -0.18%  '

$ws.Range("D5").Value = '''240.12'
$ws.Range("E5").Value = '  -0.30%  '

$ws.Range("E6").Value = '  -0.53%  '

$ws.Range("D7").Value = '''0.9987'
$ws.Range("E7").Value = '  -0.20%  '

$ws.Range("D8").Value = '''0.07595'
$ws.Range("E8").Value = '  -1.23%  '

$ws.Range("D9").Value = '''0.2900'
$ws.Range("E9").Value = '  -1.41%  '

$ws.Range("D10").Value = '''24.70'
$ws.Range("E10").Value = '  +0.78%  '

$ws.Range("D11").Value = '''0.07734'
$ws.Range("E11").Value = '  -0.16%  '

$ws.Range("D12").Value = '''5.025'
$ws.Range("E12").Value = '  +0.09%  '

$ws.Range("D13").Value = '''0.6785'
$ws.Range("E13").Value = '  -0.23%  '

$ws.Range("D14").Value = '''0.00001059'
$ws.Range("E14").Value = '  -2.33%  '

$ws.Range("D15").Value = '''82.91'
$ws.Range("E15").Value = '  -1.00%  '

$ws.Range("D16").Value = '''6.121'
$ws.Range("E16").Value = '  -0.53%  '

$ws.Range("D17").Value = '29.381.10'
$ws.Range("E17").Value = '  -0.20%  '

$ws.Range("D18").Value = '''227.63'
$ws.Range("E18").Value = '  -0.69%  '

$ws.Range("E19").Value = '  -1.05%  '

$ws.Range("D20").Value = '''0.9985'
$ws.Range("E20").Value = '  -0.21%  '

$ws.Range("D21").Value = '''7.470'
$ws.Range("E21").Value = '  +0.31%  '

$ws.Range("D22").Value = '''0.9990'
$ws.Range("E22").Value = '  -0.16%  '

$ws.Range("D23").Value = '''158.72'
$ws.Range("E23").Value = '  +0.88%  '

$ws.Range("E24").Value = '  -0.35%  '

$ws.Range("D25").Value = '''8.422'
$ws.Range("E25").Value = '  +0.51%  '

$ws.Range("D26").Value = '''17.64'
$ws.Range("E26").Value = '  -0.23%  '

$ws.Range("D27").Value = '''1.435'
$ws.Range("E27").Value = '  +9.17%  '

$ws.Range("D28").Value = '''1.457'
$ws.Range("E28").Value = '  -0.93%  '

$ws.Range("D29").Value = '''0.05593'
$ws.Range("E29").Value = '  -2.08%  '

$ws.Range("D30").Value = '''4.102'
$ws.Range("E30").Value = '  -0.32%  '

$ws.Range("E31").Value = '  +0.23%  '

$ws.Range("D32").Value = '''1.160'
$ws.Range("E32").Value = '  +0.13%  '

$ws.Range("D33").Value = '''1.829'
$ws.Range("E33").Value = '  -1.12%  '

$ws.Range("D34").Value = '''0.6954'
$ws.Range("E34").Value = '  -1.96%  '

$ws.Range("D35").Value = '''2.582'
$ws.Range("E35").Value = '  -0.26%  '

$ws.Range("D36").Value = '''0.01799'
$ws.Range("E36").Value = '  +0.01%  '

$ws.Range("D37").Value = '1.225.75'
$ws.Range("E37").Value = '  +0.08%  '

$ws.Range("D38").Value = '''2.719'
$ws.Range("E38").Value = '  -2.06%  '

$ws.Range("E39").Value = '  -1.54%  '

$ws.Range("D40").Value = '''0.8971'
$ws.Range("E40").Value = '  -1.44%  '

$ws.Range("D41").Value = '''0.9984'
$ws.Range("E41").Value = '  -0.23%  '

$ws.Range("D42").Value = '''101.38'
$ws.Range("E42").Value = '  -0.43%  '

$ws.Range("D43").Value = '''65.45'
$ws.Range("E43").Value = '  -1.22%  '

$ws.Range("D44").Value = '''7.208'
$ws.Range("E44").Value = '  +0.89%  '

$ws.Range("D45").Value = '''0.3984'
$ws.Range("E45").Value = '  -1.02%  '

$ws.Range("B46").Value = 'RenderToken'
$ws.Range("C46").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D46").Value = '''1.688'
$ws.Range("E46").Value = '  +0.07%  '

$ws.Range("B47").Value = 'EnergySwap'
$ws.Range("C47").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D47").Value = '''8.993'
$ws.Range("E47").Value = '  -0.05%  '

$ws.Range("E48").Value = '  +1.43%  '

$ws.Range("B49").Value = 'Cronos'
$ws.Range("C49").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D49").Value = '''0.05696'
$ws.Range("E49").Value = '  -0.34%  '

$ws.Range("B50").Value = 'Mantle'
$ws.Range("C50").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D50").Value = '''0.4621'
$ws.Range("E50").Value = '  -0.17%  '

$ws.Range("B51").Value = 'NEARProtocol'
$ws.Range("C51").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D51").Value = '''1.344'
$ws.Range("E51").Value = '  -0.50%  '
